$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Ready for handoff" -> "In Translation" (shared string used by Overview!E2/F2,
# zh-cn!C2 and de-de!C2 - the Status column for the row)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# Narrow the "Status" column from ~17.22 chars to ~13.41 chars on all three sheets.
# (12.5 is the ColumnWidth value that rounds to the closest on-disk width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
